# Update countries & provincias Spain
# Refreshes the COVID-19 snapshot data (columns B:H) for the countries whose
# figures moved since the previous pull, re-establishes the descending sort
# on "Casos totales" (col B) where two countries' totals crossed over
# (Moldavia/Kenia and Montserrat/Islas Malvinas), and bumps the "last
# updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4,2).Value = 5573154
$ws.Cells.Item(4,3).Value = 6522
$ws.Cells.Item(4,4).Value = 2924174
$ws.Cells.Item(4,5).Value = 2475794
$ws.Cells.Item(4,7).Value = 58
$ws.Cells.Item(4,8).Value = 173186

# --- Row 6: India ---
$ws.Cells.Item(6,2).Value = 2684314
$ws.Cells.Item(6,3).Value = 36998
$ws.Cells.Item(6,4).Value = 1939454
$ws.Cells.Item(6,5).Value = 693175
$ws.Cells.Item(6,7).Value = 640
$ws.Cells.Item(6,8).Value = 51685

# --- Row 12: Chile ---
$ws.Cells.Item(12,2).Value = 387502
$ws.Cells.Item(12,3).Value = 1556
$ws.Cells.Item(12,4).Value = 360385
$ws.Cells.Item(12,5).Value = 16604
$ws.Cells.Item(12,7).Value = 61
$ws.Cells.Item(12,8).Value = 10513

# --- Row 20: Italia ---
$ws.Cells.Item(20,2).Value = 254235
$ws.Cells.Item(20,3).Value = 320
$ws.Cells.Item(20,4).Value = 203968
$ws.Cells.Item(20,5).Value = 14867
$ws.Cells.Item(20,7).Value = 4
$ws.Cells.Item(20,8).Value = 35400

# --- Row 27: Canada ---
$ws.Cells.Item(27,2).Value = 122186
$ws.Cells.Item(27,3).Value = 99
$ws.Cells.Item(27,4).Value = 108567
$ws.Cells.Item(27,5).Value = 4593

# --- Row 46: Guatemala ---
$ws.Cells.Item(46,2).Value = 62944
$ws.Cells.Item(46,3).Value = 382
$ws.Cells.Item(46,4).Value = 51530
$ws.Cells.Item(46,5).Value = 9025
$ws.Cells.Item(46,7).Value = 10
$ws.Cells.Item(46,8).Value = 2389

# --- Rows 64/65: Moldavia overtakes Kenia in total cases, so the two rows
#     swap places to keep the sheet sorted descending by "Casos totales".
$ws.Cells.Item(64,1).Value = "Moldavia"
$ws.Cells.Item(64,2).Value = 30377
$ws.Cells.Item(64,3).Value = 194
$ws.Cells.Item(64,4).Value = 21220
$ws.Cells.Item(64,5).Value = 8249
$ws.Cells.Item(64,6).Value = 0
$ws.Cells.Item(64,7).Value = 12
$ws.Cells.Item(64,8).Value = 908

$ws.Cells.Item(65,1).Value = "Kenia"
$ws.Cells.Item(65,2).Value = 30365
$ws.Cells.Item(65,3).Value = 245
$ws.Cells.Item(65,4).Value = 17160
$ws.Cells.Item(65,5).Value = 12723
$ws.Cells.Item(65,6).Value = 0
$ws.Cells.Item(65,7).Value = 8
$ws.Cells.Item(65,8).Value = 482

# --- Row 83: Republica de Macedonia ---
$ws.Cells.Item(83,2).Value = 12840
$ws.Cells.Item(83,3).Value = 101
$ws.Cells.Item(83,5).Value = 3119
$ws.Cells.Item(83,7).Value = 3
$ws.Cells.Item(83,8).Value = 547

# --- Row 101: Grecia ---
$ws.Cells.Item(101,2).Value = 7222
$ws.Cells.Item(101,3).Value = 147
$ws.Cells.Item(101,5).Value = 3188
$ws.Cells.Item(101,7).Value = 2
$ws.Cells.Item(101,8).Value = 230

# --- Row 112: Namibia ---
$ws.Cells.Item(112,2).Value = 4344
$ws.Cells.Item(112,3).Value = 190
$ws.Cells.Item(112,4).Value = 2379
$ws.Cells.Item(112,5).Value = 1929
$ws.Cells.Item(112,7).Value = 1
$ws.Cells.Item(112,8).Value = 36

# --- Row 142: Uganda ---
$ws.Cells.Item(142,2).Value = 1560
$ws.Cells.Item(142,3).Value = 60
$ws.Cells.Item(142,4).Value = 1165
$ws.Cells.Item(142,5).Value = 380
$ws.Cells.Item(142,7).Value = 2
$ws.Cells.Item(142,8).Value = 15

# --- Row 166: Trinidad yTobago ---
$ws.Cells.Item(166,2).Value = 565
$ws.Cells.Item(166,3).Value = 13
$ws.Cells.Item(166,5).Value = 413
$ws.Cells.Item(166,7).Value = 1
$ws.Cells.Item(166,8).Value = 12

# --- Row 190: Brunei ---
$ws.Cells.Item(190,4).Value = 139
$ws.Cells.Item(190,5).Value = 0

# --- Rows 213/214: Montserrat / Islas Malvinas swap places (tie on total
#     cases resolved the other way this update).
$ws.Cells.Item(213,1).Value = "Montserrat"
$ws.Cells.Item(213,4).Value = 12
$ws.Cells.Item(213,8).Value = 1

$ws.Cells.Item(214,1).Value = "Islas Malvinas"
$ws.Cells.Item(214,4).Value = 13
$ws.Cells.Item(214,8).Value = 0

# --- Header timestamp ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 17 de Agosto de 2020 a las 17:49"
